# Update to US commit ecc67274 on 6/5/24
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOC")

# Data update: the hydro opportunity cost value in B2 changed from 70 to 15
$ws.Range("B2").Value = 15

# Reflect the author's new cell selection on the HOC sheet (was F9, now B5)
$ws.Select()
$ws.Range("B5").Select()
